$wb = $excel.ActiveWorkbook

# --- Crypto sheet: ticker list shifts (DAI-USD removed, DOT-USD added) and refreshed figures ---
$wsCrypto = $wb.Worksheets.Item("Crypto")
$wsCrypto.Range("A5").Value = "DOGE-USD"
$wsCrypto.Range("A6").Value = "DOT-USD"
$wsCrypto.Cells.Item(5, 2).Value = 0.007118002616628451
$wsCrypto.Cells.Item(5, 3).Value = 0.1316804669561276
$wsCrypto.Cells.Item(6, 2).Value = 0.001307446046721381
$wsCrypto.Cells.Item(6, 3).Value = 0.05807611652739762

# --- Stocks sheet: refreshed rendement/risque figures ---
$wsStocks = $wb.Worksheets.Item("Stocks")
$wsStocks.Cells.Item(2, 2).Value = 0.001040274805943864
$wsStocks.Cells.Item(2, 3).Value = 0.01869454577872105
$wsStocks.Cells.Item(3, 2).Value = 0.0006949433457509914
$wsStocks.Cells.Item(3, 3).Value = 0.01714595398544126
$wsStocks.Cells.Item(4, 2).Value = 0.0005167004672045829
$wsStocks.Cells.Item(4, 3).Value = 0.01551572079405033
$wsStocks.Cells.Item(5, 2).Value = 0.0007470263954237301
$wsStocks.Cells.Item(5, 3).Value = 0.01601628517350174
$wsStocks.Cells.Item(7, 2).Value = 0.0005033747232596893
$wsStocks.Cells.Item(7, 3).Value = 0.01564144545509704
$wsStocks.Cells.Item(8, 2).Value = 0.0004784300136614811
$wsStocks.Cells.Item(8, 3).Value = 0.01624235007327382
$wsStocks.Cells.Item(10, 2).Value = 0.001145376479342465
$wsStocks.Cells.Item(10, 3).Value = 0.02279879191341065
$wsStocks.Cells.Item(11, 2).Value = 0.001243009109070928
$wsStocks.Cells.Item(11, 3).Value = 0.02204348697407759
$wsStocks.Cells.Item(12, 2).Value = 0.0006150341346900171
$wsStocks.Cells.Item(12, 3).Value = 0.0159461907015981
$wsStocks.Cells.Item(13, 2).Value = 0.0005622893994287648
$wsStocks.Cells.Item(13, 3).Value = 0.02659193029580943
$wsStocks.Cells.Item(15, 2).Value = 0.0005219524393068289
$wsStocks.Cells.Item(15, 3).Value = 0.02046117811023034
$wsStocks.Cells.Item(16, 2).Value = 0.0001591414459888042
$wsStocks.Cells.Item(16, 3).Value = 0.01531359660101454
$wsStocks.Cells.Item(17, 2).Value = 0.0004404400974900765
$wsStocks.Cells.Item(17, 3).Value = 0.02080516373726509
$wsStocks.Cells.Item(18, 2).Value = 0.0001977318121623644
$wsStocks.Cells.Item(18, 3).Value = 0.02160799037043204
$wsStocks.Cells.Item(19, 2).Value = 0.0002721477072497585
$wsStocks.Cells.Item(19, 3).Value = 0.007540395733154363
$wsStocks.Cells.Item(21, 2).Value = 0.0003834017122983347
$wsStocks.Cells.Item(21, 3).Value = 0.01617251175688755
$wsStocks.Cells.Item(22, 2).Value = 0.0008541702915496731
$wsStocks.Cells.Item(22, 3).Value = 0.01398766111829008
$wsStocks.Cells.Item(24, 2).Value = 0.0005836712815873848
$wsStocks.Cells.Item(24, 3).Value = 0.01633593916632693
$wsStocks.Cells.Item(25, 2).Value = 0.000643053007192623
$wsStocks.Cells.Item(25, 3).Value = 0.01964230048833777
$wsStocks.Cells.Item(26, 2).Value = 0.0007814675614885127
$wsStocks.Cells.Item(26, 3).Value = 0.01512624906783477
$wsStocks.Cells.Item(27, 2).Value = 0.00002182021188460377
$wsStocks.Cells.Item(27, 3).Value = 0.0180290367383379
$wsStocks.Cells.Item(28, 2).Value = 0.0001522794625631943
$wsStocks.Cells.Item(28, 3).Value = 0.02294337554453492
$wsStocks.Cells.Item(31, 2).Value = 0.0002242077814322748
$wsStocks.Cells.Item(31, 3).Value = 0.0131977588250851
$wsStocks.Cells.Item(32, 2).Value = 0.0007011405672606856
$wsStocks.Cells.Item(32, 3).Value = 0.01583273226599517
$wsStocks.Cells.Item(33, 2).Value = 0.0005118191814721143
$wsStocks.Cells.Item(33, 3).Value = 0.01870428613517179
$wsStocks.Cells.Item(34, 2).Value = 0.0005017443887235558
$wsStocks.Cells.Item(34, 3).Value = 0.01497661687076971
$wsStocks.Cells.Item(35, 2).Value = 0.0002883859169647865
$wsStocks.Cells.Item(35, 3).Value = 0.01627308676744977
$wsStocks.Cells.Item(36, 2).Value = 0.0002410526056787282
$wsStocks.Cells.Item(36, 3).Value = 0.01534361462446466
$wsStocks.Cells.Item(37, 2).Value = 0.0004214885173098026
$wsStocks.Cells.Item(37, 3).Value = 0.02127118682541714
$wsStocks.Cells.Item(38, 2).Value = 0.0003987861640579937
$wsStocks.Cells.Item(38, 3).Value = 0.01173882212740428
$wsStocks.Cells.Item(39, 2).Value = 0.0006394180108161035
$wsStocks.Cells.Item(39, 3).Value = 0.01794304293207898
$wsStocks.Cells.Item(40, 2).Value = 0.0003650689361128029
$wsStocks.Cells.Item(40, 3).Value = 0.01170302295075284
$wsStocks.Cells.Item(41, 2).Value = 0.0007474156643399355
$wsStocks.Cells.Item(41, 3).Value = 0.0151853932899955
$wsStocks.Cells.Item(42, 2).Value = 0.001114912019626893
$wsStocks.Cells.Item(42, 3).Value = 0.01754759844811884
$wsStocks.Cells.Item(43, 2).Value = 0.0005876770654467698
$wsStocks.Cells.Item(43, 3).Value = 0.01458630688395485
$wsStocks.Cells.Item(44, 2).Value = 0.000874116228186846
$wsStocks.Cells.Item(44, 3).Value = 0.01785245366507381
$wsStocks.Cells.Item(45, 2).Value = 0.0006758199255548094
$wsStocks.Cells.Item(45, 3).Value = 0.0132308602556366
$wsStocks.Cells.Item(46, 2).Value = 0.0002353456120550179
$wsStocks.Cells.Item(46, 3).Value = 0.01496592226932979
$wsStocks.Cells.Item(47, 2).Value = 0.0005220993984260933
$wsStocks.Cells.Item(47, 3).Value = 0.01385985442728419
$wsStocks.Cells.Item(48, 2).Value = 0.001164256288947322
$wsStocks.Cells.Item(48, 3).Value = 0.01764232710184882
$wsStocks.Cells.Item(49, 2).Value = 0.0006105529046968781
$wsStocks.Cells.Item(49, 3).Value = 0.01511814561634551
$wsStocks.Cells.Item(51, 2).Value = 0.000486770779842394
$wsStocks.Cells.Item(51, 3).Value = 0.01883846632750532
$wsStocks.Cells.Item(52, 2).Value = 0.002645968338421937
$wsStocks.Cells.Item(52, 3).Value = 0.03113618432632288
$wsStocks.Cells.Item(53, 2).Value = 0.0008808812895971894
$wsStocks.Cells.Item(53, 3).Value = 0.01712494443486227
$wsStocks.Cells.Item(54, 2).Value = 0.0002869198445529987
$wsStocks.Cells.Item(54, 3).Value = 0.01212595722966888
$wsStocks.Cells.Item(55, 2).Value = 0.0006747962661647984
$wsStocks.Cells.Item(55, 3).Value = 0.01705886620269276
$wsStocks.Cells.Item(56, 2).Value = 0.0004697779166595275
$wsStocks.Cells.Item(56, 3).Value = 0.01206818851698195
$wsStocks.Cells.Item(57, 2).Value = 0.00027953951258137
$wsStocks.Cells.Item(57, 3).Value = 0.0147574114256138
$wsStocks.Cells.Item(58, 2).Value = 0.0004755092624498337
$wsStocks.Cells.Item(58, 3).Value = 0.01207479412072682
$wsStocks.Cells.Item(59, 2).Value = 0.0003809285866428235
$wsStocks.Cells.Item(59, 3).Value = 0.01495336095858413
$wsStocks.Cells.Item(61, 2).Value = 0.0003837763352706466
$wsStocks.Cells.Item(61, 3).Value = 0.02216581016518313
$wsStocks.Cells.Item(62, 2).Value = 0.0001987648500499688
$wsStocks.Cells.Item(62, 3).Value = 0.01346611503024537
$wsStocks.Cells.Item(63, 2).Value = 0.0002358167096391119
$wsStocks.Cells.Item(63, 3).Value = 0.01651415209573394
$wsStocks.Cells.Item(64, 2).Value = 0.0004139167453992826
$wsStocks.Cells.Item(64, 3).Value = 0.01326156426985187
$wsStocks.Cells.Item(65, 2).Value = 0.0005213024922746487
$wsStocks.Cells.Item(65, 3).Value = 0.01708457339245396
$wsStocks.Cells.Item(66, 2).Value = 0.0004730921216594786
$wsStocks.Cells.Item(66, 3).Value = 0.01715971856654771
$wsStocks.Cells.Item(67, 2).Value = 0.0003043076244917839
$wsStocks.Cells.Item(67, 3).Value = 0.01424449283468452
$wsStocks.Cells.Item(68, 2).Value = 0.0001286428637827233
$wsStocks.Cells.Item(68, 3).Value = 0.01466714787473131
$wsStocks.Cells.Item(69, 2).Value = 0.0002644469651902904
$wsStocks.Cells.Item(69, 3).Value = 0.01380701223177128
$wsStocks.Cells.Item(70, 2).Value = 0.0007958073266762181
$wsStocks.Cells.Item(70, 3).Value = 0.01612866966790688
$wsStocks.Cells.Item(72, 2).Value = 0.0009458049701675722
$wsStocks.Cells.Item(72, 3).Value = 0.01968057830413482
$wsStocks.Cells.Item(73, 2).Value = 0.0008162825466445646
$wsStocks.Cells.Item(73, 3).Value = 0.01819887291172259
$wsStocks.Cells.Item(74, 2).Value = 0.0008832852936993715
$wsStocks.Cells.Item(74, 3).Value = 0.01665940029502935
$wsStocks.Cells.Item(75, 2).Value = 0.0005928200250519067
$wsStocks.Cells.Item(75, 3).Value = 0.01682683958795737
$wsStocks.Cells.Item(76, 2).Value = 0.0007452283662661815
$wsStocks.Cells.Item(76, 3).Value = 0.01620161865235581
$wsStocks.Cells.Item(77, 2).Value = -0.00024891996758383
$wsStocks.Cells.Item(77, 3).Value = 0.01747335449581276
$wsStocks.Cells.Item(78, 2).Value = 0.00008346650527783939
$wsStocks.Cells.Item(78, 3).Value = 0.01208490838389468
$wsStocks.Cells.Item(79, 2).Value = 0.0001703375457198461
$wsStocks.Cells.Item(79, 3).Value = 0.02051378024748879
$wsStocks.Cells.Item(80, 2).Value = 0.0005577931824642184
$wsStocks.Cells.Item(80, 3).Value = 0.01367608902130595
$wsStocks.Cells.Item(81, 2).Value = 0.0005228493615443696
$wsStocks.Cells.Item(81, 3).Value = 0.01845269540396306
